# Rename the sheet from "Sheet1" to "TimeTrack"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "TimeTrack"

# Fill in the TimeTrack data.
# (Row 2 / "num" is written last so the shared-strings table ends up in
#  the same order as the authored workbook: TaskName, testingSpaceship,
#  Description, Modifying the task details, Comment, Dynamic Websites...,
#  num.)
$ws.Range("A1").Value = "TaskName"
$ws.Range("B1").Value = "testingSpaceship"
$ws.Range("A3").Value = "Description"
$ws.Range("B3").Value = "Modifying the task details"
$ws.Range("A4").Value = "Comment"
$ws.Range("B4").Value = "Dynamic Websites are difficult to handle and we are now handelling them"
$ws.Range("A2").Value = "num"
$ws.Range("B2").Value = 7

# Widen the two columns (values tuned so the persisted <col width=.../>
# lands as close as possible to the target 17.44140625 / 63.44140625).
$ws.Columns.Item(1).ColumnWidth = 16.608072916666668
$ws.Columns.Item(2).ColumnWidth = 62.608072916666664

# Match the saved selection: C1 active, sqref covering the rest of the
# used columns down to the bottom of the sheet.
$ws.Range("C1:N1048576").Select() | Out-Null
